$wb = $excel.ActiveWorkbook

# --- Sheet "6-4-22 " (2nd tab): only the view/selection changed ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Application.ActiveWindow.ScrollRow = 12
$ws2.Range("D12").Select()

# --- Sheet "7-4-22" (3rd tab): content + view changes ---
$ws3 = $wb.Worksheets.Item(3)

# Row 4
$ws3.Range("E4").Value = "Design Pattern Session-1 hr, Meeting with Client Rafi-45min, , Making changes in Acceptance Cretiria, Constraints, Assumptions. "
$ws3.Range("C4").Value = "Design Pattern Session, Meeting with Client Rafi, , Making changes in Acceptance Cretiria, Constraints, Assumptions. Training Head - Acceptance Cretiria, Constraints, Assumptions."
$ws3.Range("H4").Value = ""
$ws3.Range("H4").Font.Name = "Calibri"
$ws3.Range("H4").Font.Size = 11
$ws3.Range("H4").Font.Color = 0

# Row 5
$ws3.Range("D5").Value = "Previous UI design has been discarded"
$ws3.Range("E5").Value = "Worked on acceptance criteria and UI design - 1:30 mins,Prototype disscussion with rafi - 45 mins,team discussion about modification on prototype - 50 mins, Trainer Landpage,feedback page-1 hr"
$ws3.Range("F5").Value = "Abstract Pattern-40 mins,Softskill Session-45 mins,"
$ws3.Range("H5").Value = ""

# Row 6
$ws3.Range("E6").Value = ""
$ws3.Rows("6").RowHeight = 25.5

# Row 7
$ws3.Range("B7").Value = "Prototype discussion with team "
$ws3.Range("C7").Value = "Assumption, Acceptance Criteria and constraints for Trainee stories"
$ws3.Range("D7").Value = "Prototype design"
$ws3.Range("E7").Value = " Discussion with team members about UI prototype - 40 mins, Worked on Acceptance criteria ,assumption and constraints for Trainee - 2 hours,Meeting with Rafi - 1 hour,  "
$ws3.Range("F7").Value = "Design pattern session - 40 mins ,softskill session- 50 mins,Others - 1 hour 20 mins"

# Row 8
$ws3.Range("C8").Value = "Meeting with rafi (reviewed UI design )-1hr"
$ws3.Range("E8").Value = ""
$ws3.Rows("8").RowHeight = 25.5

# Row 10
$ws3.Range("B10").Value = ""
$ws3.Range("C10").Value = "Trainer - constraints, acceptance criteria, assumptions"
$ws3.Range("D10").Value = "Previous trainer UI design has been discarded"
$ws3.Range("D10").Font.Name = "Calibri"
$ws3.Range("D10").Font.Size = 10
$ws3.Range("D10").Font.Color = 0
$ws3.Range("E10").Value = "Worked on acceptance criteria and UI design - 1:30 mins,Prototype disscussion with rafi - 45 mins,team discussion about modification on prototype - 50 mins, Trainer - constraints, acceptance criteria, assumptions - 1 hr"
$ws3.Range("F10").Value = "Session with Rafi about abstract pattern - 40 mins"
$ws3.Range("F10").Font.Name = "Calibri"
$ws3.Range("F10").Font.Size = 10
$ws3.Range("F10").Font.Color = 0
$ws3.Rows("10").RowHeight = 94.5

# Row 12
$ws3.Range("A12").Value = "Arul "
$ws3.Range("B12").Value = "Refining trainee user stories"
$ws3.Range("C12").Value = "Assumption, Acceptance Criteria and constraints for Trainee stories"
$ws3.Range("D12").Value = "Prototype design"
$ws3.Range("E12").Value = "Meeting with client - 45mins,I've discussed about UI protyping with team members - 45mins, Done acceptance criteria for Trainee user stories - 120 mins"
$ws3.Range("F12").Value = "Softskill Session-45 mins, Design pattern session ,others 90mins"

# View/selection on sheet 3
$ws3.Application.ActiveWindow.ScrollRow = 1
$ws3.Range("E11").Select()

Write-Output "done"
